$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": zero out current-month (julio) category breakdown ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$s1Values = @{
    "M2" = 0
    "E3" = 0
    "I3" = 0
    "M3" = 0
    "D6" = 0
    "M6" = 0
    "E9" = 0
    "I9" = 0
    "M9" = 0
    "N9" = 0
    "M10" = 0
    "D12" = 0
    "M12" = 0
    "M26" = 0
    "L27" = 0
    "M27" = 0
    "L29" = 0
    "M31" = 0
    "P31" = 0
    "D34" = 0
    "E34" = 0
    "H34" = 0
    "I34" = 0
    "M44" = 0
    "M45" = 0
    "M49" = 0
    "M51" = 0
    "O51" = 0
}
foreach ($ref in $s1Values.Keys) {
    $ws1.Range($ref).Value = $s1Values[$ref]
}

$s1Text = @{
    "D58" = "0 de 56"
    "E58" = "0 de 56"
    "H58" = "0 de 56"
    "I58" = "0 de 56"
    "L58" = "0 de 56"
    "M58" = "0 de 56"
    "N58" = "0 de 56"
    "P58" = "0 de 56"
}
foreach ($ref in $s1Text.Keys) {
    $ws1.Range($ref).Value = $s1Text[$ref]
}

# --- Sheet "VENTA MENSUAL": roll months forward (abril..julio -> mayo..agosto) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$s2Text = @{
    "C1" = "mayo"
    "D1" = "junio"
    "E1" = "julio"
    "F1" = "agosto"
}
foreach ($ref in $s2Text.Keys) {
    $ws2.Range($ref).Value = $s2Text[$ref]
}

$s2Values = @{
    "C2" = 0
    "D2" = 1191.06
    "E2" = 2820.46
    "F2" = 0
    "C3" = 3316.15
    "D3" = 559.83
    "E3" = 7300.34
    "F3" = 0
    "C4" = 0
    "C6" = 0
    "E6" = 4990.88
    "F6" = 0
    "E9" = 2938.7
    "F9" = 0
    "C10" = 18249.33
    "D10" = 3793.28
    "E10" = 24572.5
    "F10" = 0
    "C12" = 1888.52
    "D12" = 0
    "E12" = 8385.37
    "F12" = 0
    "C16" = 6711.78
    "D16" = 9019.8
    "E16" = 0
    "C26" = 8756.98
    "D26" = 15779.69
    "E26" = 1980.92
    "F26" = 0
    "C27" = 3865.18
    "D27" = 714.15
    "E27" = 1212.58
    "F27" = 0
    "C29" = 787.97
    "D29" = 23.76
    "E29" = 1670.44
    "F29" = 0
    "D31" = 4575.95
    "E31" = 207.29
    "F31" = 0
    "D34" = 869.73
    "E34" = 1914.26
    "F34" = 0
    "D40" = 86.4
    "E40" = 0
    "C44" = 1167.85
    "D44" = 1053.12
    "E44" = 2128.28
    "F44" = 0
    "C45" = 737.72
    "D45" = 1696.97
    "E45" = 785.08
    "F45" = 0
    "C48" = 598.58
    "D48" = 0
    "C49" = 3470.45
    "D49" = 0
    "E49" = 2249.02
    "F49" = 0
    "C50" = 0
    "D50" = 84.9
    "E50" = 0
    "C51" = 0
    "D51" = 3819.46
    "E51" = 1516.62
    "F51" = 0
    "D55" = 25
    "E55" = 0
    "C57" = 808.39
    "D57" = 0
    "C58" = 50358.9
    "D58" = 43293.1
    "E58" = 64672.74
    "F58" = 0
}
foreach ($ref in $s2Values.Keys) {
    $ws2.Range($ref).Value = $s2Values[$ref]
}

# --- Column width changes on "VENTA MENSUAL" ---
# (ColumnWidth COM units read back ~0.8333 below the stored XML character
# width, so subtract that offset to land on the exact target XML widths
# of 14 / 13 / 12 for columns C / E / F respectively.)
$ws2.Columns.Item(3).ColumnWidth = 13.166666666666666
$ws2.Columns.Item(5).ColumnWidth = 12.166666666666666
$ws2.Columns.Item(6).ColumnWidth = 11.166666666666666